$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "collection_id" column (hidden column B, holding the Resource
# Identifier field) was removed from the Digital Object bulk-import
# template. Select the full column (as a user would click the column
# header) and delete it, shifting every later column (and the shared
# strings / cell styles they reference) one place to the left.
$ws.Range("B1:B1048576").Select()
$ws.Columns("B").Delete()
